$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.740.16"
$ws.Range("E2").Value = "  -3.43%  "
$ws.Range("D3").Value = "2.911.61"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.59%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D9").Value = "2.910.05"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("E10").Value = "  +6.04%  "
$ws.Range("E11").Value = "  -4.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.13%  "
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.64%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "3.390.84"
$ws.Range("E16").Value = "  -4.06%  "
$ws.Range("D17").Value = "60.710.00"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.29%  "
$ws.Range("D19").Value = "2.910.89"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.684"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("E26").Value = "  -4.59%  "
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.33%  "
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.95%  "
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.30%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("E41").Value = "  -5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.295"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "373.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.26%  "
$ws.Range("D47").Value = "2.666.37"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.107"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.00%  "
